$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new roster entries (Thomas / Pilon, both rank "FF") as new rows
# 83 and 84, mirroring the look of the existing roster rows above: a
# bottom-bordered, left-aligned name cell in column A, the rank in column
# B, and FALSE checkboxes across the qualification columns C:K.

$ws.Range("A82:K82").Copy($ws.Range("A83:K83")) | Out-Null
$ws.Range("A82:K82").Copy($ws.Range("A84:K84")) | Out-Null

$ws.Range("A83").Value = "Thomas"
$ws.Range("B83").Value = "FF"
$ws.Range("C83:K83").Value = $false
$ws.Range("A83").NumberFormat = "General"

$ws.Range("A84").Value = "Pilon"
$ws.Range("B84").Value = "FF"
$ws.Range("C84:K84").Value = $false
$ws.Range("A84").NumberFormat = "General"
